# Rename the embedded logo pictures in the headers/footers.
#   - Pearson logo (both footers): image2.png -> image1.png
#   - BTec logo (header, first page): image1.jpg -> image2.jpg
$d = $word.ActiveDocument
$sec = $d.Sections.First

# Both footer stories (default + first-page) carry the Pearson logo.
$footers = $sec.Footers
for ($i = 1; $i -le $footers.Count; $i++) {
    $f = $footers.Item($i)
    if ($f.Exists -and $f.Range.InlineShapes.Count -gt 0) {
        $shp = $f.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image1.png"
        }
    }
}

# The first-page header carries the BTec logo.
$headers = $sec.Headers
for ($i = 1; $i -le $headers.Count; $i++) {
    $h = $headers.Item($i)
    if ($h.Exists -and $h.Range.InlineShapes.Count -gt 0) {
        $shp = $h.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            $shp.Name = "image2.jpg"
        }
    }
}
